$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set G3 (Invalid) to 1
$ws.Range("G3").Value = 1

# Set H3:H18 (Absent) to 1 for each attendance row
for ($r = 3; $r -le 18; $r++) {
    $ws.Cells.Item($r, 8).Value = 1
}
